$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs | Cthrc1 | Fzd5 | ECs
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Cthrc1"
$ws.Cells.Item(2,3).Value = "Fzd5"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 5.449420333333333
$ws.Cells.Item(2,8).Value = 16.348261
$ws.Cells.Item(2,9).Value = 0.8985142489564721
$ws.Cells.Item(2,10).Value = 0.8985142489564723
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 10.000565
$ws.Cells.Item(2,14).Value = 30.001695
$ws.Cells.Item(2,15).Value = 0.6316353758144477
$ws.Cells.Item(2,16).Value = 0.6316353758144477
$ws.Cells.Item(2,17).Value = 54.49728225582167
$ws.Cells.Item(2,18).Value = 490.475540302395
$ws.Cells.Item(2,19).Value = 0.5675333853142576
$ws.Cells.Item(2,20).Value = 0.5675333853142576

# Row 3: FAPs | Cthrc1 | Fzd5 | FAPs
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Cthrc1"
$ws.Cells.Item(3,3).Value = "Fzd5"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 5.449420333333333
$ws.Cells.Item(3,8).Value = 16.348261
$ws.Cells.Item(3,9).Value = 0.8985142489564721
$ws.Cells.Item(3,10).Value = 0.8985142489564723
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 4.264793333333333
$ws.Cells.Item(3,14).Value = 12.79438
$ws.Cells.Item(3,15).Value = 0.2693642149089528
$ws.Cells.Item(3,16).Value = 0.2693642149089528
$ws.Cells.Item(3,17).Value = 23.24065150813111
$ws.Cells.Item(3,18).Value = 209.16586357318
$ws.Cells.Item(3,19).Value = 0.2420275852546675
$ws.Cells.Item(3,20).Value = 0.2420275852546675

# Row 4: FAPs | Cthrc1 | Fzd5 | sCs
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Cthrc1"
$ws.Cells.Item(4,3).Value = "Fzd5"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 5.449420333333333
$ws.Cells.Item(4,8).Value = 16.348261
$ws.Cells.Item(4,9).Value = 0.8985142489564721
$ws.Cells.Item(4,10).Value = 0.8985142489564723
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.567455
$ws.Cells.Item(4,14).Value = 4.702364999999999
$ws.Cells.Item(4,15).Value = 0.09900040927659938
$ws.Cells.Item(4,16).Value = 0.09900040927659938
$ws.Cells.Item(4,17).Value = 8.541721148584999
$ws.Cells.Item(4,18).Value = 76.87549033726499
$ws.Cells.Item(4,19).Value = 0.08895327838754705
$ws.Cells.Item(4,20).Value = 0.08895327838754706

# Row 5: sCs | Cthrc1 | Fzd5 | ECs
$ws.Cells.Item(5,1).Value = "sCs"
$ws.Cells.Item(5,2).Value = "Cthrc1"
$ws.Cells.Item(5,3).Value = "Fzd5"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.6155033333333333
$ws.Cells.Item(5,8).Value = 1.84651
$ws.Cells.Item(5,9).Value = 0.1014857510435278
$ws.Cells.Item(5,10).Value = 0.1014857510435278
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 10.000565
$ws.Cells.Item(5,14).Value = 30.001695
$ws.Cells.Item(5,15).Value = 0.6316353758144477
$ws.Cells.Item(5,16).Value = 0.6316353758144477
$ws.Cells.Item(5,17).Value = 6.155381092716667
$ws.Cells.Item(5,18).Value = 55.39842983444999
$ws.Cells.Item(5,19).Value = 0.06410199050019018
$ws.Cells.Item(5,20).Value = 0.06410199050019018

# Row 6: sCs | Cthrc1 | Fzd5 | FAPs
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Cthrc1"
$ws.Cells.Item(6,3).Value = "Fzd5"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.6155033333333333
$ws.Cells.Item(6,8).Value = 1.84651
$ws.Cells.Item(6,9).Value = 0.1014857510435278
$ws.Cells.Item(6,10).Value = 0.1014857510435278
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 4.264793333333333
$ws.Cells.Item(6,14).Value = 12.79438
$ws.Cells.Item(6,15).Value = 0.2693642149089528
$ws.Cells.Item(6,16).Value = 0.2693642149089528
$ws.Cells.Item(6,17).Value = 2.624994512644444
$ws.Cells.Item(6,18).Value = 23.6249506138
$ws.Cells.Item(6,19).Value = 0.02733662965428531
$ws.Cells.Item(6,20).Value = 0.02733662965428531

# Row 7: sCs | Cthrc1 | Fzd5 | sCs
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Cthrc1"
$ws.Cells.Item(7,3).Value = "Fzd5"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 0.6666666666666666
$ws.Cells.Item(7,7).Value = 0.6155033333333333
$ws.Cells.Item(7,8).Value = 1.84651
$ws.Cells.Item(7,9).Value = 0.1014857510435278
$ws.Cells.Item(7,10).Value = 0.1014857510435278
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.567455
$ws.Cells.Item(7,14).Value = 4.702364999999999
$ws.Cells.Item(7,15).Value = 0.09900040927659938
$ws.Cells.Item(7,16).Value = 0.09900040927659938
$ws.Cells.Item(7,17).Value = 0.9647737773499998
$ws.Cells.Item(7,18).Value = 8.682963996149999
$ws.Cells.Item(7,19).Value = 0.01004713088905233
$ws.Cells.Item(7,20).Value = 0.01004713088905233
